$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Clear the stray empty "INNING_NUMBER" cells on the "ODI Batting" sheet
#    (rows 6 and 8) so they are completely absent, as in the target file.
# ---------------------------------------------------------------------------
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiBatting.Range("B6").ClearContents()
$odiBatting.Range("B8").ClearContents()

# ---------------------------------------------------------------------------
# 2. Add a new worksheet "ODI Batting Extra" right after "ODI Batting".
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $odiBatting)
$newSheet.Name = "ODI Batting Extra"

# Header row
$newSheet.Cells.Item(1, 1).Value = "MATCH_CODE"
$newSheet.Cells.Item(1, 2).Value = "BATTING_POSITION"
$newSheet.Cells.Item(1, 3).Value = "NUM_4"
$newSheet.Cells.Item(1, 4).Value = "NUM_6"
$newSheet.Cells.Item(1, 5).Value = "PERCENT_RUNS_OF_TOTAL"
$newSheet.Cells.Item(1, 6).Value = "MAN_OF_MATCH"

# Reuse the same bold / bordered header style used on the other sheets
# (style index 1 in xl/styles.xml) by copying formats from an existing
# header cell rather than creating a brand new style entry.
$odiBatting.Range("A1").Copy()
$newSheet.Range("A1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column A (MATCH_CODE) holds numeric-looking codes stored as text.
$newSheet.Range("A2:A13").NumberFormat = "@"
# Columns C, D and E also hold numeric-looking values stored as text.
$newSheet.Range("C2:E13").NumberFormat = "@"

$data = @(
    @{ Row = 2;  Code = "4415"; Pos = 4;     Num4 = "2"; Num6 = "0"; Pct = "16.49%"; Mom = "NO" },
    @{ Row = 3;  Code = "4419"; Pos = $null; Num4 = "";  Num6 = "";  Pct = "";       Mom = "NO" },
    @{ Row = 4;  Code = "4421"; Pos = $null; Num4 = "";  Num6 = "";  Pct = "";       Mom = "NO" },
    @{ Row = 5;  Code = "4460"; Pos = 5;     Num4 = "3"; Num6 = "3"; Pct = "21.23%"; Mom = "NO" },
    @{ Row = 6;  Code = "4474"; Pos = $null; Num4 = "";  Num6 = "";  Pct = "";       Mom = "NO" },
    @{ Row = 7;  Code = "4475"; Pos = 5;     Num4 = "1"; Num6 = "0"; Pct = "5.26%";  Mom = "NO" },
    @{ Row = 8;  Code = "4478"; Pos = $null; Num4 = "";  Num6 = "";  Pct = "";       Mom = "NO" },
    @{ Row = 9;  Code = "4487"; Pos = 5;     Num4 = "0"; Num6 = "0"; Pct = "4.20%";  Mom = "NO" },
    @{ Row = 10; Code = "4517"; Pos = $null; Num4 = "";  Num6 = "";  Pct = "";       Mom = "NO" },
    @{ Row = 11; Code = "4550"; Pos = 1;     Num4 = "2"; Num6 = "1"; Pct = "7.61%";  Mom = "NO" },
    @{ Row = 12; Code = "4557"; Pos = 3;     Num4 = "4"; Num6 = "2"; Pct = "29.74%"; Mom = "NO" },
    @{ Row = 13; Code = "4559"; Pos = 3;     Num4 = "1"; Num6 = "0"; Pct = "5.84%";  Mom = "NO" }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $newSheet.Cells.Item($r, 1).Value = $entry.Code
    if ($entry.Pos -ne $null) {
        $newSheet.Cells.Item($r, 2).Value = $entry.Pos
    } else {
        $newSheet.Cells.Item($r, 2).NumberFormat = "@"
        $newSheet.Cells.Item($r, 2).Value = ""
    }
    $newSheet.Cells.Item($r, 3).Value = $entry.Num4
    $newSheet.Cells.Item($r, 4).Value = $entry.Num6
    $newSheet.Cells.Item($r, 5).Value = $entry.Pct
    $newSheet.Cells.Item($r, 6).Value = $entry.Mom
}

$newSheet.Range("A1").Select()
